$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) runs/balls/fours/sixes for rows 2-6 (columns C-F)
# before overwriting anything, so the reassignment below is based on the
# original data, not on values already moved by this same script.
$row2 = @($ws.Cells.Item(2, 3).Value(), $ws.Cells.Item(2, 4).Value(), $ws.Cells.Item(2, 5).Value(), $ws.Cells.Item(2, 6).Value())
$row3 = @($ws.Cells.Item(3, 3).Value(), $ws.Cells.Item(3, 4).Value(), $ws.Cells.Item(3, 5).Value(), $ws.Cells.Item(3, 6).Value())
$row4 = @($ws.Cells.Item(4, 3).Value(), $ws.Cells.Item(4, 4).Value(), $ws.Cells.Item(4, 5).Value(), $ws.Cells.Item(4, 6).Value())
$row5 = @($ws.Cells.Item(5, 3).Value(), $ws.Cells.Item(5, 4).Value(), $ws.Cells.Item(5, 5).Value(), $ws.Cells.Item(5, 6).Value())
$row6 = @($ws.Cells.Item(6, 3).Value(), $ws.Cells.Item(6, 4).Value(), $ws.Cells.Item(6, 5).Value(), $ws.Cells.Item(6, 6).Value())

function Set-InningsRow($rowIndex, $values) {
    $ws.Cells.Item($rowIndex, 3).Value = $values[0]
    $ws.Cells.Item($rowIndex, 4).Value = $values[1]
    $ws.Cells.Item($rowIndex, 5).Value = $values[2]
    $ws.Cells.Item($rowIndex, 6).Value = $values[3]
}

# Innings rows are reordered: 2<-3, 3<-4, 4<-2, 5<-6, 6<-5
Set-InningsRow 2 $row3
Set-InningsRow 3 $row4
Set-InningsRow 4 $row2
Set-InningsRow 5 $row6
Set-InningsRow 6 $row5
